# "make the excel more beautiful" - reformat the bugs sheet, refresh sample
# data and tidy up column widths / row heights / font & border styling.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Refresh the sample / placeholder data in rows 2 and 3
# ---------------------------------------------------------------------
$ws.Range("A2").Value = "/dm/browse/HSTZYJF-21170"
$ws.Range("B2").Value = "描述1"
$ws.Range("C2").ClearContents()
$ws.Range("D2").ClearContents()
$ws.Range("E2").ClearContents()
$ws.Range("F2").Value = "某人"
$ws.Range("G2").Value = "什么什么问题"

$ws.Range("A3").Value = "browse/HSTZYJF-21169"
$ws.Range("B3").Value = "分配给你"
$ws.Range("C3").Value = "解决(修复成功)"
$ws.Range("D3").Value = "添加修改方案"
$ws.Range("E3").Value = "添加修改备注"

Write-Output "values updated"
